$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text storage for numeric-looking strings in columns D and E by
# forcing Text format before assignment, then clearing the format afterwards
# so the saved style matches the original (unstyled) cells.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.692.59'
$ws.Range("E2").Value = '  +0.78%  '
$ws.Range("D3").Value = '1.850.54'
$ws.Range("E3").Value = '  +0.67%  '
$ws.Range("D4").Value = '1.035'
$ws.Range("E4").Value = '  +0.62%  '
$ws.Range("D5").Value = '322.33'
$ws.Range("E5").Value = '  +1.11%  '
$ws.Range("D6").Value = '1.031'
$ws.Range("E6").Value = '  +0.50%  '
$ws.Range("D7").Value = '0.4388'
$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("E8").Value = '  +1.48%  '
$ws.Range("D9").Value = '0.07390'
$ws.Range("E9").Value = '  +0.40%  '
$ws.Range("D10").Value = '0.8823'
$ws.Range("E10").Value = '  +1.01%  '
$ws.Range("D11").Value = '21.56'
$ws.Range("E11").Value = '  +0.48%  '
$ws.Range("D12").Value = '1.864.56'
$ws.Range("E12").Value = '  +1.51%  '
$ws.Range("D13").Value = '5.499'
$ws.Range("E13").Value = '  +0.27%  '
$ws.Range("D14").Value = '6.701'
$ws.Range("E14").Value = '  +0.47%  '
$ws.Range("D15").Value = '0.07162'
$ws.Range("E15").Value = '  +0.25%  '
$ws.Range("D16").Value = '84.98'
$ws.Range("E16").Value = '  +2.90%  '
$ws.Range("D17").Value = '1.037'
$ws.Range("E17").Value = '  +0.51%  '
$ws.Range("D18").Value = '0.000009053'
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("E19").Value = '  +0.40%  '
$ws.Range("D20").Value = '15.46'
$ws.Range("E20").Value = '  +0.53%  '
$ws.Range("D21").Value = '27.700.80'
$ws.Range("E21").Value = '  +0.77%  '
$ws.Range("E22").Value = '  +0.75%  '
$ws.Range("D23").Value = '11.31'
$ws.Range("E23").Value = '  +1.21%  '
$ws.Range("D24").Value = '2.086.70'
$ws.Range("E24").Value = '  +1.81%  '
$ws.Range("D25").Value = '2.073'
$ws.Range("E25").Value = '  +7.55%  '
$ws.Range("D26").Value = '158.86'
$ws.Range("E26").Value = '  +1.01%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("E28").Value = '  +2.89%  '
$ws.Range("D29").Value = '5.317'
$ws.Range("E29").Value = '  +1.42%  '
$ws.Range("D30").Value = '117.65'
$ws.Range("E30").Value = '  +1.44%  '
$ws.Range("D31").Value = '0.09056'
$ws.Range("E31").Value = '  -0.24%  '
$ws.Range("D32").Value = '0.7712'
$ws.Range("E32").Value = '  +0.71%  '
$ws.Range("E33").Value = '  +0.22%  '
$ws.Range("D34").Value = '3.005'
$ws.Range("E34").Value = '  +4.71%  '
$ws.Range("D35").Value = '4.556'
$ws.Range("E35").Value = '  +1.46%  '
$ws.Range("E36").Value = '  +0.32%  '
$ws.Range("D37").Value = '1.147'
$ws.Range("E37").Value = '  +0.70%  '
$ws.Range("D38").Value = '0.01973'
$ws.Range("E38").Value = '  +0.10%  '
$ws.Range("D39").Value = '0.05262'
$ws.Range("E39").Value = '  +0.26%  '
$ws.Range("D40").Value = '2.845'
$ws.Range("E40").Value = '  +2.10%  '
$ws.Range("E41").Value = '  +0.18%  '
$ws.Range("D42").Value = '0.1668'
$ws.Range("E42").Value = '  +0.16%  '
$ws.Range("D43").Value = '6.860'
$ws.Range("E43").Value = '  +3.27%  '
$ws.Range("D44").Value = '8.694'
$ws.Range("E44").Value = '  +2.32%  '
$ws.Range("E45").Value = '  +1.28%  '
$ws.Range("D46").Value = '10.72'
$ws.Range("E46").Value = '  +1.70%  '
$ws.Range("D49").Value = '1.699'
$ws.Range("E49").Value = '  -0.28%  '
$ws.Range("D50").Value = '0.4690'
$ws.Range("E50").Value = '  +1.14%  '
$ws.Range("D51").Value = '1.886'
$ws.Range("E51").Value = '  -0.34%  '

# Rows 47 and 48 swap rank (PaxDollar now ranks above Cronos) with updated data.
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = '1.033'
$ws.Range("E47").Value = '  +0.50%  '

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = '0.06588'
$ws.Range("E48").Value = '  +3.93%  '

$dataRange.ClearFormats()

